# Populate the "passport_trips" sheet with header row + trip records, then
# apply header (bold/border/centered) and date (YYYY-MM-DD) styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row values
# ---------------------------------------------------------------------
$headers = @("country", "city", "visited", "start_date", "end_date", "latitude", "longitude")
$col = 1
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# ---------------------------------------------------------------------
# 2. Data rows (country, city, visited, start_date, end_date, latitude, longitude)
# ---------------------------------------------------------------------
$rows = @(
    @("GR", "Athens",       $true,  45981, 45988, 37.9755648,  23.7348324),
    @("GR", "Athens",       $true,  45981, 45988, 37.9755648,  23.7348324),
    @("GR", "Athens",       $true,  45981, 45988, 37.9755648,  23.7348324),
    @("CO", "Bogota",       $true,  45422, 45432, 4.6533817,   -74.0836331),
    @("CO", "Bogota",       $true,  45422, 45432, 4.6533817,   -74.0836331),
    @("AR", "Buenos Aires", $false, 46082, 46090, -34.6095579, -58.3887904)
)

$r = 2
foreach ($row in $rows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Header style: bold font, thin box border, centered/top alignment.
#    Build the style on A1 alone, then fan it out to the rest of row 1
#    with a copy/paste-special so only one new style record is created.
# ---------------------------------------------------------------------
$headerCell = $ws.Range("A1")
$headerCell.Font.Bold = $true
$headerCell.HorizontalAlignment = -4108   # xlCenter
$headerCell.VerticalAlignment = -4160     # xlTop
$headerCell.Borders.LineStyle = 1         # xlContinuous (thin)
$headerCell.Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)    # xlPasteFormats

# ---------------------------------------------------------------------
# 4. Date style: custom number format "YYYY-MM-DD" for start_date/end_date.
#    Set it twice (lowercase then uppercase) on the anchor cell so both
#    format codes get registered, then fan the final style out to the
#    rest of the date columns.
# ---------------------------------------------------------------------
$dateCell = $ws.Range("D2")
$dateCell.NumberFormat = "yyyy-mm-dd"
$dateCell.NumberFormat = "YYYY-MM-DD"
$dateCell.Copy()
$ws.Range("D2:E7").PasteSpecial(-4122)    # xlPasteFormats
